# Apply the target changes to validBranches.xlsx:
#  1) Replace the numeric values in column B (rows 1-10) with their updated
#     (shortened) amounts.
#  2) Highlight the bordered data range (A1:E10) with a solid white fill
#     (new fill is added to the style table and applied to the existing
#     bordered cell styles).
#  3) Recolor the thin cell borders around the data range from white to red.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Update the numeric values in column B -----------------------------
# (last digit of each original 4-digit amount is dropped)
$newValues = @(123, 432, 876, 543, 809, 147, 987, 309, 294, 907)

for ($i = 0; $i -lt $newValues.Count; $i++) {
    $row = $i + 1
    $ws.Range("B$row").Value = $newValues[$i]
}

# --- 2) & 3) Restyle the A1:E10 data range --------------------------------
$dataRange = $ws.Range("A1:E10")

# New solid white fill applied across the whole bordered range.
$dataRange.Interior.Color = 16777215

# Existing thin borders switch from white to red.
$dataRange.Borders.Color = 255
